# 自动更新Excel文件 - 2026-02-05 23:18:15
# For every data row, decrement the "remaining days" (column E) by one.
# When the remaining days would drop to zero (i.e. current value is 1 or
# less), the cycle restarts: E is reset to the "total days" value (column D)
# and the "start date" (column F, stored as a plain yyyymmdd integer) is
# advanced by that many days.
#
# Rows whose F value is not a well-formed 8-digit yyyymmdd date are left
# untouched (mirrors source data row 36, which has a malformed date and
# was skipped by the original update).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count()
if ($lastRow -lt 2) { $lastRow = 2 }

for ($row = 2; $row -le $lastRow; $row++) {
    $dCell = $ws.Cells.Item($row, 4)   # column D - total days
    $eCell = $ws.Cells.Item($row, 5)   # column E - remaining days
    $fCell = $ws.Cells.Item($row, 6)   # column F - start date (yyyymmdd)

    $dVal = $dCell.Value()
    $eVal = $eCell.Value()
    $fVal = $fCell.Value()

    if ($null -eq $dVal -or $null -eq $eVal -or $null -eq $fVal) {
        continue
    }

    $totalDays = [int]$dVal
    $remaining = [int]$eVal
    $fText = [string][int]$fVal

    # Only touch rows with a proper 8-digit yyyymmdd date value.
    if ($fText.Length -ne 8) {
        continue
    }

    $year = [int]$fText.Substring(0, 4)
    $month = [int]$fText.Substring(4, 2)
    $day = [int]$fText.Substring(6, 2)

    try {
        $startDate = Get-Date -Year $year -Month $month -Day $day
    } catch {
        continue
    }

    if ($remaining -le 1) {
        $newRemaining = $totalDays
        $newDate = $startDate.AddDays($totalDays)
        $newFVal = [int]($newDate.ToString("yyyyMMdd"))
    } else {
        $newRemaining = $remaining - 1
        $newFVal = [int]$fText
    }

    $eCell.Value = $newRemaining
    $fCell.Value = $newFVal
}
